$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared-string values in A2 / B2
$ws.Range("A2").Value = "passive income"
$ws.Range("B2").Value = "passive.income.nadi.myfirstdrawermenuproject"

# Adjust row 2 height
$ws.Rows.Item(2).RowHeight = 23.85

# Move the active selection to A3
$ws.Range("A3").Select()
